$wb = $excel.ActiveWorkbook

# --- Logged Week 17 data: append new game totals to the running per-game logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "6 8 0 -2 15 2 2 -2 5 3 6 0 2 5 12 2 3 3 0 4 2 2 0 16 2 5 -2 13 -1 4 5 1 11 3 10 1 4 5 1 3 2 4 14 0 5 1 3 4 14 4 3 0 6 1 7 15 22 3 1 2 2 44 5 3 30 5 6 11 3 6 4 18 4 8 11 1 4 5 24 3 4 1 8 1 2 5 1 1 5 1 50 1 2 7 2 3 1 3 4 0 13 8 6 9 11 -2 0 14 5 42 3 4 1 1 3 -4 4 0 6 5 3 34 0 3 10 3 14 3 20 4 3 3 0 0 1 -2 9 4 -2 2 3 7 2 7 5 -2 0 10 16 5 6 2 1 -2 24 -1 4 6 1 37 0 5 -4 3 1 10 5 21 4 4 -1 1 4 3 4 3 28 1 6 5 25 6 12 3 14 3 5 1 4 2 0 4 8 9 4 8 2 4 1 4 2 0 3 -2 3 7 4 9 5 20 15 3 2 2 3 -1 3 -2 7 -1 1 8 -2 5 2 8 4 4 3 3 4 5 4 1 0 1 5 5 1 1 3 8 5 8 1 9 4 5 -1 3 5 -2 -1 4 4 3 1 5 7 5 3 3 9 5 1 3 5 6 1 2 1 8 5 -2 7 11 0 8 3 2 5 0 2 10 7 1 2 6 9 -3 1 11 4 7 2 2 1 3 2 3 2 3 2 5 0 1 8 8 14 2 1 -4 1 5 6 11 3 1 -1 -1 39 4 4 1 19 3 0 6 0 5 7 -2 -2 6 0 1 9 3 10 7 5 18 0 1 36 4 7 2 0 2 30 -1 4 2 2 7 7 5 5 14 21 9 -2 4 6 5 0 7 0 5 24 3 5 5 7 4 13 6 1 5 11 4 10 13 1 20 11 2 4 17 1 44 19 3 1 2 4 1 4 0 5 8 6 2 -3 3 4 11 7 1 8 2 9 3 5 2 3 -3 5 0 7 6 3 0 4 3 18 5 4 0 7 2 2 6 19 1 17 5 4 3 5 6 2 7 8 5 1 18 3 6 2 6 17 2 9 5 7 3 2 7 7 20 32 0 7 1 25 5 6 3 6 0 3 0 7 0 0 4 3 1 13 7 2 22 2 12 3 5 6 3 6 5 5 5 8 22 9 0 3 2 5 3 2 3 18 6 27 20 4 4 7 4 4 11 4 8 72 16 6 1 16 2 2 5 2 0 1 2 0 1 3 2 4 0 3 3 5 0 4 8 1 15 2 48 8 10 -1 9 23 2 4 4 -2 1 3 12 2 14 4 33 5 4 5 12 10 9 -1 5 7 1 2 5 0 3 3 0 1 2 13 -4 9 5 2 15 1 1 6 -4 0 19 3 9 2 11 -2 3 18 2 2 3 -3 0 4 1 16 3 2 4 11 6 -1 15 11 4 2 2 5 5 8 1 0 28 2 4 2 3 2 6 5 12 31 9 7 1 9 9 8 0 0 21 3 4 5 20 14 -4 3 4 -2 20 6 2 8 5 9 2 1 11 5 9 -1 3 1 1 2 2 -2 2 5 0 31 6 2 0 4 8 8 1 9 0 5 7 8 5 12 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 0 1 3 3 2 -2 13 4 3 -1 9 2 4 0 12 5 -1 3 6 4 -3 6 6 6 2 8 6 16 6 1 11 14 3 3 22 2 3 5 5 7 3 6 1 5 4 1 7 2 2 0 -3 5 7 3 0 9 -1 9 3 2 12 2 10 -5 4 -2 -1 12 3 10 3 7 1 16 2 1 7 7 2 10 7 2 5 6 5 8 5 6 4 4 2 6 3 4 9 3 5 6 7 11 6 12 0 1 7 1 4 9 4 6 1 3 8 7 8 5 10 1 6 5 4 2 3 19 5 12 1 4 19 4 12 2 8 4 2 -2 0 0 5 2 11 0 -4 1 3 1 11 5 7 5 10 1 -1 -3 4 1 4 2 5 4 2 19 0 3 8 2 3 4 5 2 5 4 6 5 0 2 1 6 -1 2 2 3 1 2 4 1 1 6 1 0 3 4 0 13 9 5 4 2 0 5 1 0 4 2 4 1 4 3 2 8 13 1 3 11 8 3 1 4 6 5 4 1 12 5 18 -1 3 5 3 8 -2 4 5 3 8 6 8 0 -1 -2 13 3 2 2 3 1 5 1 4 5 4 3 -4 -1 7 2 18 4 8 5 0 11 13 21 5 3 2 1 1 2 4 1 7 2 -2 5 3 16 4 6 8 12 6 13 4 -1 6 10 7 3 15 2 8 8 2 6 1 0 7 -2 3 -5 3 4 2 2 -7 1 2 1 2 -3 5 3 7 1 1 1 9 5 17 11 3 -1 12 2 4 3 3 32 6 6 4 6 3 2 3 -2"
$ws.Range("C2").Value = "11 0 4 2 29 22 1 3 8 3 0 0 0 4 1 3 2 10 10 9 4 3 4 3 1 4 0 1 1 8 2 5 2 5 6 9 1 -2 4 2 2 22 1 3 4 5 1 1 -1 9 4 16 1 7 3 11 2 0 -1 0 12 6 5 3 3 1 2 4 2 2 5 1 10 2 2 -1 3 7 -1 3 9 1 -1 2 9 5 3 7 3 2 0 1 5 5 1 6 3 1 -3 2 2 2 2 6 3 1 2 -1 4 3 7 2 0 7 2 0 2 2 -1 0 2 4 2 5 1 20 3 1 26 5 74 5 -1 40 5 1 -1 1 9 4 4 0 18 6 -2 3 3 1 -3 1 1 3 0 11 3 4 0 4 16 2 1 9 12 3 1 2 4 3 11 5 4 2 1 14 7 4 7 3 0 6 7 4 3 4 4 12 5 1 13 4 16 25 4 12 3 4 1 -2 -1 0 9 1 1 0 -1 8 8 1 2 -2 4 -5 5 2 4 8 8 3 0 21 3 3 3 3 0 2 5 1 0 2 -2 24 1 11 3 1 12 4 3 8 4 4 29 5 2 1 3 3 -2 7 3 1 2 3 10 3 3 2 13 7 2 1 1 7 2 9 2 13 4 1 2 10 3 -1 0 10 14 3 2 2 3 2 3 3 0 -1 0 3 12 3 3 7 2 2 7 4 2 3 8 14 1 11 5 4 5 -1 2 10 2 4 5 7 3 6 13 1 2 2 10 3 0 6 -1 2 1 0 -1 2 0 2 0 1 2 2 11 2 -1 17 3 5 6 1 1 9 3 3 3 15 3 2 1 3 4 13 1 11 2 -2 55 2 -2 4 6 15 9 2 4 17 4 2 6 1 1 -3 5 4 1 1 0 4 2 4 1 3 2 0 8 3 0 4 3 9 -1 12 6 3 -2 -1 0 3 1 31 -1 2 8 0 13 0 0 0 2 2 2 4 1 15 3 -2 6 0 1 9 1 7 -1 1 3 9 4 2 7 6 15 1 3 -5 4 -2 -2 11 1 4 14 2 6 4 14 4 0 9 -2 2 1 2 3 12 4 0 4 0 1 1 -1 0 1 -4 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 3 11 4 4 3 13 2 2 -1 2 2 19 13 0 4 4 8 22 -1 7 2 1 -1 3 -4 1 -1 4 2 5 2 10 -2 2 1 3 -2 2 9 0 -1 -1 9 0 0 5 1 0 0 10 7 21 3 46 3 -3 -5 2 2 1 11 15 -3 4 66 1 3 3 0 0 1 2 -1 1 0 -1 0 -8 1 24 1 3 -1 3 1 -1 5 4 2 4 2 1 -1 4 1 14 0 11 0 1 7 2 1 8 9 2 11 5 -6 3 10 7 5 6 4 5 5 -1 6 11 2 0 -3 4 1 2 7 2 2 -1 3 4 -1 5 1 9 3 -3 5 2 7 0 5 -1 5 4 2 3 6 2 2 0 1 8 3 2 0 8 13 3 5 8 -2 2 1 3 3 2 5 1 6 15 1 2 13 -2 5 6 3 4 11 8 11 1 2 7 1 0 -2 1 4 2 2 4 2 3 7 9 8 -1 4 11 5 4 9 3 5 5 6 0 1 0 5 8 -2 19 3 1 1 9 1 1 1 3 12 2 2 2 8 2 2 2 6 2 4 1 1 2 3 14 1 7 7 -2 11 4 7 1 8 1 -4"
$ws.Range("B3").Value = "12 47 9 25 21 4 16 20 13 16 11 9 4 21 19 12 -1 2 7 10 9 6 -3 4 29 14 4 18 11 1 19 25 11 7 13 5 13 22 5 7 0 0 7 4 6 13 7 7 8 5 19 6 6 4 10 1 33 31 3 25 13 9 18 22 4 4 10 10 7 16 9 9 14 5 14 2 7 2 9 12 9 12 3 7 5 8 30 11 17 7 9 32 9 11 5 18 2 -1 20 10 17 10 9 14 34 6 3 11 6 10 18 39 24 3 8 4 6 4 5 4 7 13 9 -1 14 14 9 21 20 12 4 11 9 5 17 13 21 19 6 11 31 2 3 8 9 9 12 6 7 4 14 8 6 18 10 11 3 1 5 14 30 5 4 6 6 31 14 5 3 6 8 15 16 9 9 5 -1 5 6 7 3 8 7 70 1 3 4 12 8 4 38 4 2 20 9 2 4 10 39 6 15 7 13 14 14 6 4 11 10 18 28 4 27 -5 11 2 11 6 6 34 -1 44 12 3 -1 8 13 13 9 6 11 6 11 20 0 7 25 6 10 4 8 27 7 5 10 43 5 9 18 5 2 9 8 7 17 5 28 9 19 -6 11 4 6 9 3 17 20 9 17 4 1 11 31 11 21 8 30 5 5 7 14 1 8 9 5 2 4 9 11 29 2 8 0 21 6 6 29 29 10 13 4 10 6 4 13 6 49 9 10 11 17 -4 7 20 8 9 7 11 4 4 12 8 27 20 42 19 9 10 6 5 16 11 6 21 37 19 29 3 9 41 7 24 5 18 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 5 5 17 15 28 7 28 2 6 11 12 8 8 7 17 3 19 7 43 27 18 7 6 15 5 17 10 9 20 8 4 14 4 11 7 7 5 7 4 4 8 16 8 9 21 1 10 6 4 11 9 8 20 13 6 1 24 6 23 19 20 18 7 1 35 39 25 3 13 13 -1 10 12 7 11 13 13 15 7 -1 9 4 5 12 22 3 3 8 8 18 11 5 22 6 8 7 7 20 20 7 13 7 4 11 15 7 5 6 13 3 11 9 9 11 0 8 7 8 4 8 20 5 8 30 6 9 15 11 1 16 5 9 8 -2 11 11 2 3 12 9 15 2 9 8 -2 3 9 4 12 7 21 6 29 8 2 7 5 6 9 8 13 8 8 4 7 39 13 11 5 -1 0 10 3 14 16 18 0 9 9 8 2 10 29 3 10 15 15 3 17 6 5 22 18 13 5 6 11 2 -1 5 3 5 4 9 4 20 10 32 4 8 -3 20 36 12 3 11 2 6 5 6 4 8 14 3 43 5 5 4 22 8 4 8 7 9 2 7 1 9 12 5 5 3 5 5 5 3 3 6 12 3 11 17 5 15 6 4 9 -1 16 12 10 -1 18 28 11 11 4 7 12 17 25 18 14 4 6 7 14 5 15 13 5 1 14 9 18 5 33 2 14 15 4 15 7 11 4 5"
$ws.Range("C3").Value = "5 6 16 6 1 8 13 9 1 16 28 8 5 21 3 -2 4 0 9 12 4 7 13 6 6 13 38 10 6 5 12 9 13 15 9 6 6 15 4 24 11 7 11 23 3 11 33 15 11 13 15 2 5 7 29 5 18 -4 20 10 49 9 5 7 24 7 6 9 7 6 18 13 11 2 17 5 9 17 -2 3 11 8 18 11 1 13 -1 40 4 3 12 0 7 4 4 19 4 13 8 8 11 11 9 16 4 39 14 26 12 -2 14 6 7 10 4 16 8 11 3 5 12 9 9 7 -6 9 8 3 19 13 8 6 4 50 6 11 3 5 9 10 10 10 1 16 10 0 7 24 18 14 9 8 13 7 14 0 9 6 4 5 2 11 8 6 2 3 13 20 13 17 5 12 10 17 5 9 1 5 -2 3 18 6 6 5 6 4 21 3 4 19 20 7 4 4 24 26 8 6 12 1 4 14 7 8 1 -1 7 6 2 12 50 5 12 16 9 5 9 22 25 14 14 11 11 13 1 4 5 3 17 27 6 11 5 4 7 22 8 13 -1 5 11 4 1 5 6 10 8 7 8 1 7 7 5 8 5 3 3 3 16 4 9 -4 0 2 14 13 11 13 12 13 6 6 15 9 4 6 13 13 9 11 6 4 19 19 12 8 12 -4 20 10 12 4 11 9 4 4 11 37 7 18 1 9 14 9 26 6 5 12 8 10 30 7 16 22 8 8 -2 12 10 2 3 12 11 21 8 6 17 22 2 28 10 8 9 4 13 10 7 14 13 9 11 11 4 15 16 19 20 8 -5 7 0 9 16 13 11 3 12 14 15 1 8 21 8 7 3 28 4 14 5 9 10 35 5 10 6 7 18 12 3 -2 10 6 5 13 4 16 1 11 21 2 11 10 2 4 19 6 14 20 9 3 4 7 14 10 1 4 9 13 6 5 6 12 9 24 8 15 13 7 4 12 16 2 3 4 21 4 5 16 9 7 9 37 10 10 20 18 11 27 32 31 6 7 19 10 33 17 11 14 10 4 14 20 14 5 11 40 5 46 3 8 3 23 13 7 11 3 4 5 8 1 15 10 5 12 11 14 19 9 2 -3 11 9 6 24 19 22 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 76 11 15 6 8 15 11 7 -2 5 9 16 -1 31 42 21 24 18 4 4 29 23 2 15 11 11 8 3 18 7 5 8 26 1 16 1 5 10 21 15 -1 5 4 5 10 6 11 23 6 5 9 55 7 4 6 13 26 7 17 21 32 7 82 25 2 11 11 27 6 50 5 3 7 7 20 4 18 5 9 5 6 27 1 2 16 9 11 7 21 2 6 19 2 52 11 9 1 6 14 35 23 6 64 6 29 22 22 23 60 10 2 3 12 9 4 23 11 9 36 10 0 -1 41 11 7 4 16 6 38 8 20 16 8 3 18 7 12 4 17 5 14 7 7 4 9 40 29 25 9 2 5 9 5 5 -1 8 3 9 4 34 5 5 9 7 6 27 9 1 12 5 7 13 17 3 3 4 17 6 11 31 7 22 10 3 13 3 24 6 1 16 9 7 6 9 25 14 11 23 -6 10 18 20 3 11 16 5 17 3 68 20 17 9 11 15 3 52 1 19 3 18 9 3 2 5 14 13 0 12 4 22 3 10 20 4 52 10 18 7 19 6 4 9 22 11 17 11 18 2 35 12 9 2 8 3 21 15 24 4 15 7"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "65 64 62 60 61 61 65 67 55 63 68 70 63 63 59 65 62 58 61 59 68 65 62 55 62 65 61 59 65 61 66 65 57 63 65 57 65 55 48 67 58 63 64 65 61 63 60 67 59 67 63 61 39 65 61 64 58 63 63 47 40 65 63 65 53 56 59 59 43 59"
$ws.Range("B5").Value = "17 19 10 20 22 21 27 22 12 24 30 37 14 14 19 20 26 6 15 66 29 22 16 19 24 23 19 18 18 16 25 21 26 13 20 8 33 3 9 25 17 19 26 20 16 18 16 36 10 23 98 17 0 15 14 19 13 21 0 13 6 35 21 20 0 9 25 20 0 18"
$ws.Range("B6").Value = "38 26 21 20 16 26 93 30 37 20 26 24 37 17 26 18 28 26 24 21 5 16 23 20 23 23 21 23 27 23 47 21 24 24 21 25 25 26 35 0 19 12 23 24 26 34 21 22 15 24 0 22"
$ws.Range("D3").Value = "49 37 39 53 41 48 49 50 49 41 51 44 42 42 49 57 59 45 46 26 43 49 36 44 39 50 43 40 44 43 48 44 58 49 40 40 46 42 58 46 35 52 46 42 34 46 50 35 41 28 40 55 45 53 43 43 48 23 37 48 55 47 42 49 41 49 46 54 48 35 34 52 55 46 49 57 42 49 53 47 49 53 33 48 49 40 31 50 46 39 35 51 45 37 53 43 39 24 41 39 37 42 44 35 57 45 53 33 50 55 45 46 57 60 44 44 42 47 43 43 36 42 61"
$ws.Range("D4").Value = "1 0 0 10 0 0 7 0 -3 0 0 4 0 0 0 11 5 0 0 0 0 4 5 0 0 11 9 0 0 -2 3 0 0 0 0 0 0 0 0 0 0 6 0 0 0 0 5 0 0 0 6 28 -2 6 0 0 2 0 0 8 0 3 0 0 7 0 13 11 0 0 0 0 5 8 0 9 0 6 0 13 11 0 0 13 7 0 0 0 0 0 0 0 6 0 0 0 0 0 4 0 0 0 0 0 4 0 0 0 -3 14 5 6 8 8 1 6 0 6 10 0 0 0 18"
$ws.Range("D5").Value = "0 17 9 0 0 8 16 0 10 5 13 0 0 0 0 10 1 -1 0 10 0 3 5 10 11 9 0 0 0 0 15 0 0 0 9 4 0 0 0 0 0 7 0 11 0 0 0 0 0 0 20 32 17 19 0 14 0 0 0 6 0 0 0 0 0 0 0 0 0 0 0 0 22 15 30 0 0 16 0 0 0 0 13 29 0 12 0 0 10 0 42 0 0 0 0 0 7 0 14 13 17 0 0 16 21 0 0 0 0 14 0 0 0 0 0 9 0 0 0 7 12 0 0 0 13 22 0 0 1 0 0 0 12 0 0 0 0 0 0 0 0 0 0 0 0"

# --- Updated season totals/splits after folding in Week 17 ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 474
$ws.Range("D2").Value = 27
$ws.Range("F2").Value = 139
$ws.Range("G2").Value = 149
$ws.Range("J2").Value = 116
$ws.Range("L2").Value = 498
$ws.Range("M2").Value = 323
$ws.Range("O2").Value = 46
$ws.Range("P2").Value = 29
$ws.Range("Q2").Value = 1071
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 343
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 228
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 59
$ws.Range("I3").Value = 118
$ws.Range("J3").Value = 70
$ws.Range("N3").Value = 37

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 390
$ws.Range("F2").Value = 95
$ws.Range("G2").Value = 99
$ws.Range("J2").Value = 41
$ws.Range("L2").Value = 563
$ws.Range("M2").Value = 368
$ws.Range("O2").Value = 53
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 963
$ws.Range("B3").Value = 23
$ws.Range("C3").Value = 381
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 238
$ws.Range("G3").Value = 71
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 121
$ws.Range("J3").Value = 109
$ws.Range("N3").Value = 37

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 185
$ws.Range("D2").Value = 123
$ws.Range("F2").Value = 385
$ws.Range("G2").Value = 381
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 177
$ws.Range("L2").Value = 110
$ws.Range("M2").Value = 99
$ws.Range("B3").Value = 115

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 9
$ws.Range("E2").Value = 18
$ws.Range("D3").Value = 19

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 46
